# Update country COVID-19 stats table to "25 de Marzo de 2020 a las 14:46" snapshot
# (title date bump, refreshed case counts, and reshuffled country-name column
#  caused by the underlying data being re-sorted by total cases).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Marzo de 2020 a las 14:46"
$ws.Cells.Item(6, 2).Value = 55081
$ws.Cells.Item(6, 3).Value = 225
$ws.Cells.Item(6, 5).Value = 53917
$ws.Cells.Item(8, 2).Value = 35353
$ws.Cells.Item(8, 3).Value = 2362
$ws.Cells.Item(8, 5).Value = 31632
$ws.Cells.Item(8, 7).Value = 22
$ws.Cells.Item(8, 8).Value = 181
$ws.Cells.Item(13, 5).Value = 7659
$ws.Cells.Item(13, 7).Value = 11
$ws.Cells.Item(13, 8).Value = 433
$ws.Cells.Item(14, 1).Value = "Austria"
$ws.Cells.Item(14, 3).Value = 277
$ws.Cells.Item(14, 4).Value = 9
$ws.Cells.Item(14, 5).Value = 5521
$ws.Cells.Item(14, 6).Value = 28
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 8).Value = 30
$ws.Cells.Item(15, 1).Value = "Paises Bajos"
$ws.Cells.Item(15, 2).Value = 5560
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(15, 5).Value = 5282
$ws.Cells.Item(15, 6).Value = 435
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 276
$ws.Cells.Item(18, 2).Value = 2971
$ws.Cells.Item(18, 3).Value = 105
$ws.Cells.Item(18, 5).Value = 2951
$ws.Cells.Item(20, 1).Value = "Suecia"
$ws.Cells.Item(20, 2).Value = 2526
$ws.Cells.Item(20, 3).Value = 227
$ws.Cells.Item(20, 4).Value = 16
$ws.Cells.Item(20, 5).Value = 2468
$ws.Cells.Item(20, 6).Value = 158
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(20, 8).Value = 42
$ws.Cells.Item(21, 1).Value = "Australia"
$ws.Cells.Item(21, 2).Value = 2431
$ws.Cells.Item(21, 3).Value = 114
$ws.Cells.Item(21, 4).Value = 118
$ws.Cells.Item(21, 5).Value = 2304
$ws.Cells.Item(21, 6).Value = 11
$ws.Cells.Item(21, 8).Value = 9
$ws.Cells.Item(24, 2).Value = 1874
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(24, 5).Value = 1830
$ws.Cells.Item(30, 1).Value = "Chile"
$ws.Cells.Item(30, 2).Value = 1142
$ws.Cells.Item(30, 3).Value = 220
$ws.Cells.Item(30, 4).Value = 22
$ws.Cells.Item(30, 5).Value = 1118
$ws.Cells.Item(30, 6).Value = 7
$ws.Cells.Item(30, 8).Value = 2
$ws.Cells.Item(31, 1).Value = "Luxemburgo"
$ws.Cells.Item(31, 2).Value = 1099
$ws.Cells.Item(31, 4).Value = 6
$ws.Cells.Item(31, 5).Value = 1085
$ws.Cells.Item(31, 6).Value = 3
$ws.Cells.Item(31, 8).Value = 8
$ws.Cells.Item(32, 1).Value = "Ecuador"
$ws.Cells.Item(32, 2).Value = 1082
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 3
$ws.Cells.Item(32, 5).Value = 1052
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 27
$ws.Cells.Item(33, 1).Value = "Pakistan"
$ws.Cells.Item(33, 2).Value = 1022
$ws.Cells.Item(33, 3).Value = 50
$ws.Cells.Item(33, 4).Value = 21
$ws.Cells.Item(33, 5).Value = 993
$ws.Cells.Item(33, 6).Value = 5
$ws.Cells.Item(33, 7).Value = 1
$ws.Cells.Item(33, 8).Value = 8
$ws.Cells.Item(34, 1).Value = "Polonia"
$ws.Cells.Item(34, 2).Value = 957
$ws.Cells.Item(34, 3).Value = 56
$ws.Cells.Item(34, 4).Value = 2
$ws.Cells.Item(34, 5).Value = 942
$ws.Cells.Item(34, 6).Value = 3
$ws.Cells.Item(34, 7).Value = 3
$ws.Cells.Item(34, 8).Value = 13
$ws.Cells.Item(35, 1).Value = "Tailandia"
$ws.Cells.Item(35, 2).Value = 934
$ws.Cells.Item(35, 3).Value = 107
$ws.Cells.Item(35, 4).Value = 70
$ws.Cells.Item(35, 5).Value = 860
$ws.Cells.Item(35, 6).Value = 11
$ws.Cells.Item(35, 8).Value = 4
$ws.Cells.Item(38, 2).Value = 880
$ws.Cells.Item(38, 3).Value = 88
$ws.Cells.Item(38, 5).Value = 867
$ws.Cells.Item(38, 6).Value = 22
$ws.Cells.Item(41, 1).Value = "Islandia"
$ws.Cells.Item(41, 2).Value = 737
$ws.Cells.Item(41, 3).Value = 89
$ws.Cells.Item(41, 4).Value = 56
$ws.Cells.Item(41, 5).Value = 679
$ws.Cells.Item(41, 6).Value = 11
$ws.Cells.Item(41, 8).Value = 2
$ws.Cells.Item(42, 1).Value = "Crucero"
$ws.Cells.Item(42, 2).Value = 712
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 587
$ws.Cells.Item(42, 5).Value = 115
$ws.Cells.Item(42, 6).Value = 15
$ws.Cells.Item(42, 8).Value = 10
$ws.Cells.Item(43, 1).Value = "Sudafrica"
$ws.Cells.Item(43, 2).Value = 709
$ws.Cells.Item(43, 3).Value = 155
$ws.Cells.Item(43, 4).Value = 12
$ws.Cells.Item(43, 5).Value = 697
$ws.Cells.Item(43, 6).Value = 2
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(44, 1).Value = "Rusia"
$ws.Cells.Item(44, 2).Value = 658
$ws.Cells.Item(44, 3).Value = 163
$ws.Cells.Item(44, 4).Value = 29
$ws.Cells.Item(44, 5).Value = 628
$ws.Cells.Item(44, 6).Value = 8
$ws.Cells.Item(44, 8).Value = 1
$ws.Cells.Item(46, 2).Value = 606
$ws.Cells.Item(46, 3).Value = 70
$ws.Cells.Item(46, 4).Value = 42
$ws.Cells.Item(46, 5).Value = 554
$ws.Cells.Item(51, 5).Value = 341
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 21
$ws.Cells.Item(58, 5).Value = 328
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 7
$ws.Cells.Item(60, 1).Value = "Irak"
$ws.Cells.Item(60, 2).Value = 346
$ws.Cells.Item(60, 3).Value = 30
$ws.Cells.Item(60, 4).Value = 89
$ws.Cells.Item(60, 5).Value = 230
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 8).Value = 27
$ws.Cells.Item(61, 1).Value = "Libano"
$ws.Cells.Item(61, 3).Value = 15
$ws.Cells.Item(61, 4).Value = 8
$ws.Cells.Item(61, 5).Value = 321
$ws.Cells.Item(61, 6).Value = 4
$ws.Cells.Item(61, 8).Value = 4
$ws.Cells.Item(62, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(62, 2).Value = 333
$ws.Cells.Item(62, 3).Value = 85
$ws.Cells.Item(62, 4).Value = 52
$ws.Cells.Item(62, 5).Value = 279
$ws.Cells.Item(62, 6).Value = 2
$ws.Cells.Item(62, 8).Value = 2
$ws.Cells.Item(114, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(114, 3).Value = 3
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 8).Value = 2
$ws.Cells.Item(116, 1).Value = "Cuba"
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 1
$ws.Cells.Item(116, 6).Value = 2
$ws.Cells.Item(116, 8).Value = 1
$ws.Cells.Item(123, 1).Value = "Honduras"
$ws.Cells.Item(123, 3).Value = 6
$ws.Cells.Item(124, 1).Value = "Mayotte"
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(128, 1).Value = "Kenia"
$ws.Cells.Item(129, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(141, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(141, 3).Value = 4
$ws.Cells.Item(142, 1).Value = "Uganda"
$ws.Cells.Item(142, 3).Value = 5
$ws.Cells.Item(143, 4).Value = 8
$ws.Cells.Item(143, 5).Value = 5
$ws.Cells.Item(152, 1).Value = "Dominica"
$ws.Cells.Item(154, 1).Value = "Seychelles"
$ws.Cells.Item(155, 1).Value = "Haiti"
$ws.Cells.Item(158, 1).Value = "Benin"
$ws.Cells.Item(159, 1).Value = "Bermudas"
$ws.Cells.Item(160, 1).Value = "Islas Caimanes"
$ws.Cells.Item(161, 1).Value = "Curazao"
$ws.Cells.Item(167, 1).Value = "Santa Sede"
$ws.Cells.Item(168, 1).Value = "Guinea"
$ws.Cells.Item(169, 1).Value = "Congo"
$ws.Cells.Item(170, 1).Value = "Suazilandia"
$ws.Cells.Item(172, 1).Value = "Liberia"
$ws.Cells.Item(174, 1).Value = "Mozambique"
$ws.Cells.Item(175, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(177, 1).Value = "Laos"
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(178, 1).Value = "San Bartolome"
$ws.Cells.Item(179, 1).Value = "Republica del Chad"
$ws.Cells.Item(180, 1).Value = "Birmania"
$ws.Cells.Item(180, 3).Value = 0
$ws.Cells.Item(181, 1).Value = "Angola"
$ws.Cells.Item(182, 1).Value = "Nepal"
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 8).Value = 0
$ws.Cells.Item(183, 1).Value = "Zimbabue"
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 8).Value = 1
$ws.Cells.Item(184, 1).Value = "Sudan"
$ws.Cells.Item(185, 1).Value = "Gambia"
$ws.Cells.Item(187, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(187, 3).Value = 2
$ws.Cells.Item(188, 1).Value = "Nicaragua"
$ws.Cells.Item(189, 1).Value = "Mali"
$ws.Cells.Item(190, 1).Value = "Mauritania"
$ws.Cells.Item(190, 3).Value = 0
$ws.Cells.Item(191, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(192, 1).Value = "Eritrea"
$ws.Cells.Item(193, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(194, 1).Value = "Siria"
$ws.Cells.Item(195, 1).Value = "Montserrat"
$ws.Cells.Item(196, 1).Value = "Timor Oriental"
$ws.Cells.Item(197, 1).Value = "Somalia"
$ws.Cells.Item(198, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(199, 1).Value = "Belice"
$ws.Cells.Item(201, 1).Value = "Granada"
$ws.Cells.Item(202, 1).Value = "Libia"

Write-Host "Applied 224 cell updates"
